# Update column headers (row 1) with the new capitalized Spanish labels
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Rama"
$ws.Range("B1").Value = "VAB"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Rama código"
$ws.Range("F1").Value = "Año"
$ws.Range("G1").Value = "Rama descripción"

# Update the URI / code row (row 2)
$ws.Range("A2").Value = "iaest-measure:rama"
$ws.Range("B2").Value = "iaest-measure:vab"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refPeriod"
$ws.Range("G2").Value = "iaest-measure:rama-descripcion"

# Update the dim/medida role row (row 3)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "medida"

# Update the datatype row (row 4)
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:double"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "xsd:date"
$ws.Range("G4").Value = "xsd:string"

# Move the mapping-ano.xlsx reference from A5 (ano column) to F5 (now the Año column)
$ws.Range("A5").Delete(-4162)
$ws.Range("F5").Value = "mapping-ano.xlsx"

# Reapply the standard cell style to the relocated cell so it matches its neighbours
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
